$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 9167
$ws.Range("I18").Value = 10600
$ws.Range("K18").Value = 10600
$ws.Range("M18").Value = -10316
$ws.Range("H40").Value = 6294.067
$ws.Range("H53").Value = 67034.47
$ws.Range("I53").Value = 715.2727
$ws.Range("J53").Value = 105429.79
$ws.Range("K53").Value = 715.2727
$ws.Range("L53").Value = 105429.79
$ws.Range("M53").Value = -78.27269999999999
$ws.Range("N53").Value = -106703.79
$ws.Range("H64").Value = 8432
$ws.Range("I64").Value = 3530.6667
$ws.Range("J64").Value = 13333.333
$ws.Range("K64").Value = 3530.6667
$ws.Range("L64").Value = 13333.333
$ws.Range("M64").Value = -3282.6667
$ws.Range("N64").Value = -13829.333
$ws.Range("H67").Value = 8432
$ws.Range("I67").Value = 3530.6667
$ws.Range("J67").Value = 13333.333
$ws.Range("K67").Value = 3530.6667
$ws.Range("L67").Value = 13333.333
$ws.Range("M67").Value = -2672.6667
$ws.Range("N67").Value = -15049.333
$ws.Range("H86").Value = 10973.5
$ws.Range("I86").Value = 9997.5
$ws.Range("J86").Value = 11949.5
$ws.Range("K86").Value = 9997.5
$ws.Range("L86").Value = 11949.5
$ws.Range("M86").Value = -8874.5
$ws.Range("N86").Value = -14195.5
$ws.Range("H89").Value = 10973.5
$ws.Range("I89").Value = 9997.5
$ws.Range("J89").Value = 11949.5
$ws.Range("K89").Value = 49987.5
$ws.Range("L89").Value = 59747.5
$ws.Range("M89").Value = -44371.5
$ws.Range("N89").Value = -70979.5
$ws.Range("H111").Value = 266.33334
$ws.Range("I111").Value = 266.33334
$ws.Range("K111").Value = 799.0000200000001
$ws.Range("M111").Value = 2267.99998
$ws.Range("H132").Value = 1270.7428
$ws.Range("I132").Value = 1048.7742
$ws.Range("K132").Value = 3146.3226
$ws.Range("M132").Value = -616.3226000000004
$ws.Range("H138").Value = 3815.6428
$ws.Range("J138").Value = 4439.4
$ws.Range("L138").Value = 13318.2
$ws.Range("N138").Value = -23598.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2348.3298
$ws.Range("I32").Value = 2201.6667
$ws.Range("K32").Value = 2201.6667
$ws.Range("M32").Value = -1914.6667
$ws.Range("H61").Value = 5056334
$ws.Range("I61").Value = 6949293
$ws.Range("K61").Value = 6949293
$ws.Range("M61").Value = -6949081
$ws.Range("H134").Value = 93833
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").Value = $null
$ws.Range("H136").Value = 5056334
$ws.Range("I136").Value = 6949293
$ws.Range("K136").Value = 20847879
$ws.Range("M136").Value = -20845329
$ws.Range("H140").Value = 99998.5
$ws.Range("J140").Value = 99998.5
$ws.Range("L140").Value = 99998.5
$ws.Range("N140").Value = -110358.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4370.7188
$ws.Range("I31").Value = 3097.9565
$ws.Range("K31").Value = 3097.9565
$ws.Range("M31").Value = -2802.9565
$ws.Range("H34").Value = 4370.7188
$ws.Range("I34").Value = 3097.9565
$ws.Range("K34").Value = 3097.9565
$ws.Range("M34").Value = -2895.9565
$ws.Range("H41").Value = 16612.857
$ws.Range("I41").Value = 16525
$ws.Range("K41").Value = 16525
$ws.Range("M41").Value = -16097
$ws.Range("H48").Value = 40000
$ws.Range("J48").Value = 40000
$ws.Range("L48").Value = 40000
$ws.Range("N48").Value = -40952
$ws.Range("H58").Value = 7465.364
$ws.Range("I58").Value = 5749.077
$ws.Range("K58").Value = 5749.077
$ws.Range("M58").Value = -5546.077
$ws.Range("H136").Value = 7465.364
$ws.Range("I136").Value = 5749.077
$ws.Range("K136").Value = 17247.231
$ws.Range("M136").Value = -14697.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 150.375
$ws.Range("I14").Value = 150.375
$ws.Range("K14").Value = 451.125
$ws.Range("M14").Value = -278.125
$ws.Range("H68").Value = 2290.2856
$ws.Range("J68").Value = 3349
$ws.Range("L68").Value = 10047
$ws.Range("N68").Value = -11669
$ws.Range("H69").Value = 833.3333
$ws.Range("I69").Value = 833.3333
$ws.Range("K69").Value = 2499.9999
$ws.Range("M69").Value = -1688.9999
$ws.Range("H71").Value = 2290.2856
$ws.Range("J71").Value = 3349
$ws.Range("L71").Value = 30141
$ws.Range("N71").Value = -38253
$ws.Range("H72").Value = 833.3333
$ws.Range("I72").Value = 833.3333
$ws.Range("K72").Value = 7499.9997
$ws.Range("M72").Value = -3443.9997
$ws.Range("H81").Value = 11506.333
$ws.Range("I81").Value = 3759.5
$ws.Range("J81").Value = 27000
$ws.Range("K81").Value = 11278.5
$ws.Range("L81").Value = 81000
$ws.Range("M81").Value = -10155.5
$ws.Range("N81").Value = -83246
$ws.Range("H84").Value = 11506.333
$ws.Range("I84").Value = 3759.5
$ws.Range("J84").Value = 27000
$ws.Range("K84").Value = 33835.5
$ws.Range("L84").Value = 243000
$ws.Range("M84").Value = -28219.5
$ws.Range("N84").Value = -254232
$ws.Range("H109").Value = 78538.69500000001
$ws.Range("I109").Value = 1089.8
$ws.Range("J109").Value = 126944.25
$ws.Range("K109").Value = 3269.4
$ws.Range("L109").Value = 380832.75
$ws.Range("M109").Value = -2229.4
$ws.Range("N109").Value = -382912.75
$ws.Range("H133").Value = 873.5
$ws.Range("I133").Value = 873.5
$ws.Range("K133").Value = 2620.5
$ws.Range("M133").Value = 2439.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14210.5
$ws.Range("J70").Value = 16179.4
$ws.Range("L70").Value = 16179.4
$ws.Range("N70").Value = -16719.4
$ws.Range("H73").Value = 14210.5
$ws.Range("J73").Value = 16179.4
$ws.Range("L73").Value = 16179.4
$ws.Range("N73").Value = -18051.4
$ws.Range("H102").Value = 4958.6665
$ws.Range("I102").Value = 4045.818
$ws.Range("K102").Value = 4045.818
$ws.Range("M102").Value = -2423.818
$ws.Range("H126").Value = 3738.7144
$ws.Range("I126").Value = 3464.7693
$ws.Range("K126").Value = 10394.3079
$ws.Range("M126").Value = -7924.3079
$ws.Range("H132").Value = 3015.8286
$ws.Range("I132").Value = 2788.5
$ws.Range("K132").Value = 8365.5
$ws.Range("M132").Value = -5835.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 9999
$ws.Range("J26").Value = 9999
$ws.Range("L26").Value = 9999
$ws.Range("N26").Value = -10589
$ws.Range("H31").Value = 3498
$ws.Range("I31").Value = 2000
$ws.Range("J31").Value = 3872.5
$ws.Range("K31").Value = 2000
$ws.Range("L31").Value = 3872.5
$ws.Range("M31").Value = -1752
$ws.Range("N31").Value = -4368.5
$ws.Range("H137").Value = 81514.836
$ws.Range("I137").Value = 80797.8
$ws.Range("K137").Value = 80797.8
$ws.Range("M137").Value = -75697.8
$ws.Range("H138").Value = 80000
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 14999.5
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 14999.5
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 14999.5
$ws.Range("M12").Value = $null
$ws.Range("N12").Value = -15283.5
$ws.Range("H113").Value = 653.9
$ws.Range("I113").Value = 652.5625
$ws.Range("K113").Value = 1957.6875
$ws.Range("M113").Value = 212.3125
$ws.Range("I122").Value = 3142.1428
$ws.Range("J122").Value = 5856.6
$ws.Range("K122").Value = 9426.428400000001
$ws.Range("L122").Value = 17569.8
$ws.Range("M122").Value = -6976.428400000001
$ws.Range("N122").Value = -22469.8
$ws.Range("H132").Value = 3440.8
$ws.Range("I132").Value = 2530.6216
$ws.Range("K132").Value = 7591.864799999999
$ws.Range("M132").Value = -5061.864799999999
